# This change adds a new, more recent weekly price record for
# "Ciboulette" (row 189) to the daily logic subset sheet. Inserting the
# row pushes every subsequent record down by one row (the former last
# row, 283, becomes row 284), which matches the diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 189; this shifts rows 189..283 down to
# 190..284 and carries the existing per-column formatting along with
# them (in particular the date style on column D).
$ws.Rows.Item(189).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(189, 1).Value = 6
$ws.Cells.Item(189, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(189, 3).Value = "Metropolitana"
$ws.Cells.Item(189, 4).Value = 44460
$ws.Cells.Item(189, 5).Value = 13
$ws.Cells.Item(189, 6).Value = 100112039
$ws.Cells.Item(189, 7).Value = "Ciboulette"
$ws.Cells.Item(189, 8).Value = "Sin especificar"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 830
$ws.Cells.Item(189, 11).Value = 1000
$ws.Cells.Item(189, 12).Value = 1200
$ws.Cells.Item(189, 13).Value = 1106
$ws.Cells.Item(189, 14).Value = "`$/docena de atados"
$ws.Cells.Item(189, 15).Value = "Región Metropolitana"
$ws.Cells.Item(189, 16).Value = 369
$ws.Cells.Item(189, 17).Value = 3
$ws.Cells.Item(189, 18).Value = "Hortaliza"
